# VBA-style RGB() helper (PowerShell has no built-in RGB cmdlet); COM
# colour properties take the packed 0x00BBGGRR long used throughout the
# PowerPoint object model.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# Swap the deck's theme palette: ppt/theme/theme2.xml (the presentation's
# one reachable Theme - it backs the slide master / every slide) held the
# "Integral" scheme; ppt/theme/theme1.xml (Notes Master only) held the
# default "Office Theme" scheme. Re-point the live scheme to the Office
# Theme's twelve colours so the main deck now carries that palette.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
